$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 33
$ws.Range("F4").Value = 71
$ws.Range("F5").Value = 100
$ws.Range("F6").Value = 940
$ws.Range("F7").Value = 514
$ws.Range("F8").Value = 4893
$ws.Range("F9").Value = 4893
$ws.Range("F10").Value = 114
$ws.Range("F11").Value = 130
$ws.Range("F15").Value = 141
$ws.Range("F16").Value = 7910
$ws.Range("F17").Value = 7910
$ws.Range("F18").Value = 261
$ws.Range("F20").Value = 562
$ws.Range("F21").Value = 2061
$ws.Range("F22").Value = 6305
$ws.Range("F23").Value = 2274
$ws.Range("F28").Value = 6278
$ws.Range("F29").Value = 166
$ws.Range("F30").Value = 50
$ws.Range("F34").Value = 6651
$ws.Range("F35").Value = 3
$ws.Range("F38").Value = 5
$ws.Range("F39").Value = 7
$ws.Range("F40").Value = 0
$ws.Range("F42").Value = 30
$ws.Range("F48").Value = 468
$ws.Range("F49").Value = 2182
$ws.Range("F50").Value = 58

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 246
$ws.Range("F4").Value = 48
$ws.Range("F16").Value = 1

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 33
$ws.Range("F5").Value = 71
$ws.Range("F6").Value = 246
$ws.Range("F7").Value = 100
$ws.Range("F8").Value = 48
$ws.Range("F9").Value = 514
$ws.Range("F10").Value = 4893
$ws.Range("F11").Value = 4893
$ws.Range("F12").Value = 114
$ws.Range("F13").Value = 130
$ws.Range("F16").Value = 141
$ws.Range("F17").Value = 7910
$ws.Range("F18").Value = 7910
$ws.Range("F19").Value = 261
$ws.Range("F21").Value = 562
$ws.Range("F22").Value = 2063
$ws.Range("F25").Value = 6305
$ws.Range("F26").Value = 2274
$ws.Range("F31").Value = 6278
$ws.Range("F32").Value = 166
$ws.Range("F33").Value = 50
$ws.Range("F37").Value = 6651
$ws.Range("F41").Value = 30
$ws.Range("F47").Value = 468
$ws.Range("F49").Value = 2182
$ws.Range("F50").Value = 58
